$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status text changes (shared string used by B3/C3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: Status (C3) + Error Detail (L3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handback transform failed"
$wsZh.Range("L3").Value = "Handback file name: h5rz52tu.akt is different with handoff file name: 83857edd-a801-4b19-8109-284b4b5759f4.b5f98ee212790f516ec90bf4707cc659a3e4833b.zh-cn."

# --- de-de sheet: Status (C3) + Error Detail (L3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handback transform failed"
$wsDe.Range("L3").Value = "Handback file name: h5rz52tu.akt is different with handoff file name: 83857edd-a801-4b19-8109-284b4b5759f4.b5f98ee212790f516ec90bf4707cc659a3e4833b.de-de."
